# "opravy v ip setting"
$wb = $excel.ActiveWorkbook

# --- ip_address_list ---
$wsIpList = $wb.Worksheets.Item("ip_address_list")
$wsIpList.Range("D1").Value = "poznvv"
$wsIpList.Range("D2").Value = "poznggv`ndf`ndf`ndf"
$wsIpList.Range("E2").Value = $true

# --- ip_adress_fav_list ---
$wsIpFav = $wb.Worksheets.Item("ip_adress_fav_list")
$wsIpFav.Range("D1").Value = "poznggv`ndf`ndf`ndf"

# --- Settings ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value = 2
$wsSettings.Range("B6").Value = 1
$wsSettings.Range("A7").Value = "editovatelné(1)/ needitovatelné(0) poznámky (default)"
$wsSettings.Range("B7").Value = 1
